$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.96'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.03'
$ws.Range("D3").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06011'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.390'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8103'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9270'
$ws.Range("D8").Style = "Normal"

$ws.Range("B9").Value = 'One'

$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.01119'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '8OneONEBestin24h'

$ws.Range("B10").Value = 'WazirX'

$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1424'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '9WazirXWRX'

$ws.Range("B11").Value = 'MandalaExchangeToken'

$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07428'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03368'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B13").Value = 'BitrueCoin'

$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03032'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '12BitrueCoinBTR'

$ws.Range("B14").Value = 'BitMartToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09353'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '13BitMartTokenBMX'

$ws.Range("B15").Value = 'MCDex'

$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.951'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '14MCDexMCB'

$ws.Range("B16").Value = 'BitForexToken'

$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001598'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '15BitForexTokenBF'

$ws.Range("B17").Value = 'CoinExToken'

$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04836'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '16CoinExTokenCET'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005430'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004152'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009814'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.00007103'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.652'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.435'
$ws.Range("D23").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1295'
$ws.Range("D26").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03970'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006403'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1072'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002901'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.006723'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '43LocalTradersLCT'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005204'
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0005803'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.010'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002314'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '48BOLOBOLOWorstin24h'
